# Time sheet update: fill in the row-5 time entry (9/30/2024, 11:00 PM -> 12:00 AM)
# which had been left blank, and move the active selection off the old "last
# used" cell (F31) onto D9, matching where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: Date / Time In / Time Out for the new entry.
$ws.Range("A5").Value = 45565                  # 9/30/2024
$ws.Range("B5").Value = 0.95833333333333337    # 11:00 PM
$ws.Range("C5").Value = 1                      # 12:00 AM (midnight, next day)

# Move the selection to D9, where the user left off.
$ws.Range("D9").Select()
